# WIP: Einbau Listbox-namen sammeln
# Fill in three new timesheet rows (64-66) on the "MA ZD" sheet: finish
# the previously-empty placeholder row 64, replace the old test row 65
# with fresh data, and append a brand-new row 66 - all part of
# collecting sample employee/customer names for the new listbox.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MA ZD")

# Row 65 previously held different sample data; start it clean (also
# drops its row-level custom formatting) and drop the leftover empty
# "C" placeholder cell entirely.
$ws.Rows.Item(65).ClearFormats()
$ws.Cells.Item(65, 3).ClearContents()

# --- write the new text/values in on-sheet reading order so the ------
# --- shared-string table grows in the same order a person typing ----
# --- the rows left-to-right, row-by-row would have produced. --------
$ws.Cells.Item(64, 1).Value = 45990
$ws.Cells.Item(64, 2).Value = 60269
$ws.Cells.Item(64, 4).Value = "Bender"
$ws.Cells.Item(64, 5).Value = "Laber fasel"

$ws.Cells.Item(65, 1).Value = 45990
$ws.Cells.Item(65, 2).Value = 10100
$ws.Cells.Item(65, 4).Value = "Wolfgang"
$ws.Cells.Item(65, 5).Value = "Test Wolle"

$ws.Cells.Item(64, 8).Value = 2
$ws.Cells.Item(64, 9).Value = "C3A61789-050C-48C7-8D18-97841FAC8470"

$ws.Cells.Item(65, 8).Value = 3
$ws.Cells.Item(65, 9).Value = "4085D7BC-22F3-4F66-84DA-557C816963A7"

$ws.Cells.Item(66, 1).Value = 46023
$ws.Cells.Item(66, 2).Value = 99887
$ws.Cells.Item(66, 5).Value = "Bla"
$ws.Cells.Item(66, 4).Value = "ANewTestCustomer"
$ws.Cells.Item(66, 8).Value = 1
$ws.Cells.Item(66, 9).Value = "B887354A-BF91-494F-969D-A1D67EA3ECB2"

# --- formatting ---------------------------------------------------------
# Row 64 already carried the sheet's normal row-level format (Calibri),
# so its non-date cells already match; only the date cell needs the
# dedicated date style used throughout column A above it.
$ws.Range("A2").Copy()
$ws.Cells.Item(64, 1).PasteSpecial(-4122)   # xlPasteFormats -> style 7 (date, Calibri)

# Rows 65/66 were typed with plain/default formatting: the date cells
# get a plain date number format (default font), everything else gets
# the sheet's usual Calibri style copied from an existing data cell,
# and the last ("Zeilen-ID") column is left completely unformatted.
$ws.Cells.Item(65, 1).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(66, 1).NumberFormat = "m/d/yyyy"

$ws.Range("B2").Copy()
$ws.Cells.Item(65, 2).PasteSpecial(-4122)
$ws.Cells.Item(66, 2).PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Cells.Item(65, 4).PasteSpecial(-4122)
$ws.Cells.Item(66, 4).PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Cells.Item(65, 5).PasteSpecial(-4122)
$ws.Cells.Item(66, 5).PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Cells.Item(65, 8).PasteSpecial(-4122)
$ws.Cells.Item(66, 8).PasteSpecial(-4122)

# --- view state: scroll down towards the new rows and select C65, -----
# --- matching where the author's cursor ended up. ----------------------
$ws.Activate() | Out-Null
$ws.Range("A52").Select() | Out-Null
$ws.Range("C65").Select() | Out-Null

Write-Output "Rows 64-66 populated"
